# Updated cryptos list - applies Price (D) and Volume(1h) (E) changes per row
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "80.376.34"
$ws.Range("E2").Value = "  +4.89%  "

$ws.Range("D3").Value = "3.183.25"
$ws.Range("E3").Value = "  +1.40%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").ClearFormats()

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "210.01"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +4.17%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "626.05"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +0.09%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.276"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +27.66%  "

$ws.Range("E8").Value = "  -0.07%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.587"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +5.37%  "

$ws.Range("D10").Value = "3.183.60"
$ws.Range("E10").Value = "  +1.44%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.594"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +27.02%  "

$ws.Range("E12").Value = "  +28.37%  "

$ws.Range("E13").Value = "  +1.29%  "

$ws.Range("D14").Value = "3.765.65"
$ws.Range("E14").Value = "  +1.32%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.25"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -0.77%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "31.87"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +6.96%  "

$ws.Range("D17").Value = "80.262.86"
$ws.Range("E17").Value = "  +4.75%  "

$ws.Range("D18").Value = "3.178.84"
$ws.Range("E18").Value = "  +1.50%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "14.19"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +3.21%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "3.01"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +9.17%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "9.15"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -1.29%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "437.04"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +10.07%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.19"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +13.74%  "

$ws.Range("E24").Value = "  +6.15%  "

$ws.Range("D25").Value = "3.345.43"
$ws.Range("E25").Value = "  +1.65%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "75.92"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +3.64%  "

$ws.Range("E27").Value = "  +1.11%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.87"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +4.30%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.999"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -0.11%  "

$ws.Range("E30").Value = "  +8.14%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.998"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +0.22%  "

$ws.Range("E32").Value = "  +4.67%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "557.70"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +8.35%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.45"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -0.52%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.149"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +13.11%  "

$ws.Range("E36").Value = "  +1.78%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "22.97"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +4.58%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.123"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +19.00%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.999"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -0.04%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.407"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +5.99%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "20.74"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +3.34%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "163.02"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -0.27%  "

$ws.Range("E43").Value = "  -0.02%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.61"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +4.33%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "189.10"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -4.06%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.80"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +5.36%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.69"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +7.34%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.779"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -3.71%  "

$ws.Range("E49").Value = "  +0.66%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "42.91"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +2.87%  "

$ws.Range("E51").Value = "  +5.79%  "
